$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '26.723.83'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +0.25%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.601.53'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +0.25%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.01'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.16%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '211.64'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +0.06%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.514'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -0.24%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.00'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +0.15%  '
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +0.13%  '
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +0.14%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '19.70'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +0.65%  '
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +1.08%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.825.97'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +0.28%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.601.47'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +0.30%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.05'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +0.55%  '
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -0.21%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '65.14'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.0₃0740'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +0.30%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '210.42'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +0.21%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '1.01'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +0.18%  '
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +2.43%  '
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -0.23%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '2.27'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -3.13%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.00'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +0.00%  '
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -1.02%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.01'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +0.21%  '
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -0.63%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '15.38'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +0.55%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.0511'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -0.90%  '
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +0.04%  '
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +0.80%  '
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +0.96%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.291.82'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +0.60%  '
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +0.49%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.604'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -2.51%  '
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +10.07%  '
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -0.34%  '
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -0.44%  '
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -1.83%  '
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +0.17%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.784'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -0.03%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '62.70'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -1.59%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.737.84'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +0.13%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '90.57'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -0.14%  '
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -1.67%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.102'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -0.18%  '
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +1.53%  '
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +0.32%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '7.43'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +0.26%  '
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +0.96%  '
